{"js": "// Update the Oreo filling confidence-interval answer (Problem 13) from\n// (2.808 , 2.988) to (2.535 , 3.165) in the Solutions table, and make\n// the table's first row repeat as a header row (matching the canonical\n// export of this worksheet's table).\n\n// 1) Replace the numeric values in the answer text.\n// Both occurrences of each number need to change, so loop over every\n// match returned by search() instead of assuming a single hit.\nconst oldLow = \"2.808\";\nconst newLow = \"2.535\";\nconst oldHigh = \"2.988\";\nconst newHigh = \"3.165\";\n\nconst lowResults = context.document.body.search(oldLow, { matchCase: true });\nlowResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < lowResults.items.length; i++) {\n  lowResults.items[i].insertText(newLow, Word.InsertLocation.replace);\n}\nawait context.sync();\n\nconst highResults = context.document.body.search(oldHigh, { matchCase: true });\nhighResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < highResults.items.length; i++) {\n  highResults.items[i].insertText(newHigh, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Make the first row of the Solutions table repeat as a header row.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  table.headerRowCount = 1;\n  await context.sync();\n}\n", "ps1": "# Update the Oreo filling confidence-interval answer (Problem 13) from\n# (2.808 , 2.988) to (2.535 , 3.165) in the Solutions table, and bring\n# the table's formatting in line with the canonical export: repeat the\n# first row as a header row and normalize the table's preferred width.\n\n$d = $word.ActiveDocument\n\n# 1) Replace the numeric values in the answer text. Both occurrences of\n# each number in the document need to change (the sentence states the\n# interval twice), so Find/Replace All (Wrap = 1) is used for each.\n$find1 = $d.Content.Find\n$find1.Text = \"2.808\"\n$find1.Replacement.Text = \"2.535\"\n$find1.Execute([ref]$find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find1.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \"2.988\"\n$find2.Replacement.Text = \"3.165\"\n$find2.Execute([ref]$find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find2.Replacement.Text, 2)\n\n# 2) Make the first row of the Solutions table repeat as a header row,\n# and re-apply the table's preferred width as a percentage so it is\n# stored as the clean integer 5000 (= 100%, in fiftieths-of-a-percent\n# units) instead of the original \"5000.0\".\n$t = $d.Tables.Item(1)\n$t.Rows.Item(1).HeadingFormat = 1\n$t.PreferredWidthType = 2\n$t.PreferredWidth = 250\n"}
